$wb = $excel.ActiveWorkbook

# Rename the second sheet ("Include from Code de régulati" -> "Include from Codes pour carac")
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Include from Codes pour carac"

# Update the Metadata sheet values
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B5").Value = "Caractérisation de la population Value Set"
$ws1.Range("B8").Value = "2023-10-19T15:25:12+00:00"
$ws1.Range("B11").Value = "Caractérisation de la population Value Set."
